$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(11, 15).Value = 0
$ws.Cells.Item(12, 15).Value = 0.001764297485351562
$ws.Cells.Item(14, 15).Value = 0.1447784900665283
$ws.Cells.Item(15, 15).Value = 0
$ws.Cells.Item(16, 15).Value = 0.007411956787109375
$ws.Cells.Item(19, 15).Value = 0
$ws.Cells.Item(20, 15).Value = 0
$ws.Cells.Item(23, 15).Value = 0.01564240455627441
$ws.Cells.Item(25, 15).Value = 0.007860183715820312
$ws.Cells.Item(26, 15).Value = 0.0655670166015625
$ws.Cells.Item(28, 15).Value = 0.0004754066467285156
$ws.Cells.Item(31, 15).Value = 0
$ws.Cells.Item(32, 15).Value = 0.002124309539794922
$ws.Cells.Item(33, 15).Value = 0.01660704612731934
$ws.Cells.Item(34, 15).Value = 0.09069466590881348
$ws.Cells.Item(35, 15).Value = 0.0324099063873291
$ws.Cells.Item(37, 15).Value = 0.3126041889190674
$ws.Cells.Item(39, 15).Value = 0.01857423782348633
$ws.Cells.Item(41, 15).Value = 0.001100778579711914
$ws.Cells.Item(43, 15).Value = 0.001159906387329102
$ws.Cells.Item(47, 15).Value = 0.009741306304931641
$ws.Cells.Item(48, 15).Value = 0
$ws.Cells.Item(49, 15).Value = 0.008355617523193359
$ws.Cells.Item(50, 15).Value = 0.0019989013671875
$ws.Cells.Item(52, 15).Value = 0.002103090286254883
$ws.Cells.Item(54, 15).Value = 0.001930713653564453
$ws.Cells.Item(56, 15).Value = 0.001018047332763672
$ws.Cells.Item(58, 15).Value = 0.001003503799438477
$ws.Cells.Item(59, 15).Value = 0.007788896560668945
$ws.Cells.Item(60, 15).Value = 0.01012277603149414
$ws.Cells.Item(61, 15).Value = 0
$ws.Cells.Item(62, 15).Value = 0.002000808715820312
$ws.Cells.Item(64, 15).Value = 0.01649737358093262
$ws.Cells.Item(65, 15).Value = 0.4899814128875732
$ws.Cells.Item(69, 15).Value = 0.03379964828491211
$ws.Cells.Item(71, 15).Value = 1.781439065933228
$ws.Cells.Item(72, 15).Value = 18.85905051231384
$ws.Cells.Item(73, 15).Value = 0.03202486038208008
$ws.Cells.Item(75, 15).Value = 0.001008510589599609
$ws.Cells.Item(78, 15).Value = 0
$ws.Cells.Item(79, 15).Value = 0.01014518737792969
$ws.Cells.Item(80, 15).Value = 0.01148843765258789
$ws.Cells.Item(81, 15).Value = 0.08442187309265137
$ws.Cells.Item(84, 15).Value = 3.714205265045166
$ws.Cells.Item(85, 15).Value = 0
$ws.Cells.Item(86, 15).Value = 0.007059812545776367
$ws.Cells.Item(87, 15).Value = 0.04146647453308105
$ws.Cells.Item(88, 15).Value = 0.07651352882385254
$ws.Cells.Item(89, 15).Value = 0
$ws.Cells.Item(90, 15).Value = 0.06658291816711426
$ws.Cells.Item(91, 15).Value = 0.03330111503601074
$ws.Cells.Item(92, 15).Value = 0.02409076690673828
$ws.Cells.Item(93, 15).Value = 0.02819657325744629
$ws.Cells.Item(94, 15).Value = 0.04870009422302246
$ws.Cells.Item(95, 15).Value = 0.04923152923583984
$ws.Cells.Item(96, 15).Value = 0.04862427711486816
$ws.Cells.Item(98, 15).Value = 0.0005524158477783203
$ws.Cells.Item(99, 15).Value = 0.01656961441040039
$ws.Cells.Item(101, 15).Value = 0.01662325859069824
$ws.Cells.Item(102, 15).Value = 0.04859185218811035
$ws.Cells.Item(104, 15).Value = 0.7758736610412598
$ws.Cells.Item(106, 15).Value = 0.01524996757507324
$ws.Cells.Item(107, 15).Value = 0.017486572265625
$ws.Cells.Item(108, 15).Value = 0
$ws.Cells.Item(109, 15).Value = 0.01500630378723145
$ws.Cells.Item(110, 15).Value = 0
$ws.Cells.Item(111, 15).Value = 0.04938173294067383
$ws.Cells.Item(112, 15).Value = 0.08317947387695312
$ws.Cells.Item(113, 15).Value = 0.08281683921813965
$ws.Cells.Item(114, 15).Value = 0.07310652732849121
$ws.Cells.Item(115, 15).Value = 0.0005521774291992188
$ws.Cells.Item(116, 15).Value = 0.02829623222351074
$ws.Cells.Item(118, 15).Value = 0.01581764221191406
$ws.Cells.Item(120, 15).Value = 0.001006841659545898
$ws.Cells.Item(122, 15).Value = 0.005987644195556641
$ws.Cells.Item(123, 15).Value = 0
$ws.Cells.Item(125, 15).Value = 0.01511192321777344
$ws.Cells.Item(127, 15).Value = 0.001000642776489258
$ws.Cells.Item(134, 15).Value = 0.006601095199584961
$ws.Cells.Item(136, 15).Value = 0
$ws.Cells.Item(138, 15).Value = 0.0009911060333251953
$ws.Cells.Item(140, 15).Value = 0
$ws.Cells.Item(143, 15).Value = 0.001001596450805664
$ws.Cells.Item(145, 15).Value = 0.001076698303222656
$ws.Cells.Item(147, 15).Value = 0
$ws.Cells.Item(148, 15).Value = 0.01520085334777832
$ws.Cells.Item(149, 15).Value = 0.001610517501831055
$ws.Cells.Item(150, 15).Value = 0.2167835235595703
$ws.Cells.Item(151, 15).Value = 0
$ws.Cells.Item(152, 15).Value = 0
$ws.Cells.Item(162, 15).Value = 0
$ws.Cells.Item(163, 15).Value = 0
$ws.Cells.Item(164, 15).Value = 0.0009837150573730469
$ws.Cells.Item(165, 15).Value = 0
$ws.Cells.Item(168, 15).Value = 0
$ws.Cells.Item(170, 15).Value = 0.007110357284545898
$ws.Cells.Item(171, 15).Value = 0.0222017765045166
$ws.Cells.Item(173, 15).Value = 0.004760265350341797
$ws.Cells.Item(174, 15).Value = 0.005360603332519531
$ws.Cells.Item(176, 15).Value = 0.00394749641418457
